$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 50
$ws.Range("F3").Value = 7121
$ws.Range("F4").Value = 3478
$ws.Range("F6").Value = 3804
$ws.Range("F7").Value = 60
$ws.Range("F9").Value = 66
$ws.Range("F11").Value = 124
$ws.Range("F14").Value = 111
$ws.Range("F15").Value = 350
$ws.Range("F16").Value = 11
$ws.Range("F19").Value = 4055
$ws.Range("F21").Value = 400
$ws.Range("F22").Value = 1023
$ws.Range("F23").Value = 526
$ws.Range("F24").Value = 1613
$ws.Range("F25").Value = 105
$ws.Range("F26").Value = 89
$ws.Range("F27").Value = 2953
$ws.Range("F28").Value = 2132
$ws.Range("F29").Value = 56
$ws.Range("F30").Value = 75
$ws.Range("F31").Value = 85
$ws.Range("F32").Value = 39
$ws.Range("F33").Value = 36
$ws.Range("F35").Value = 4153
$ws.Range("F36").Value = 428
$ws.Range("F37").Value = 314
$ws.Range("F38").Value = 54
$ws.Range("F39").Value = 938
$ws.Range("F40").Value = 739
$ws.Range("F41").Value = 149
$ws.Range("F43").Value = 1584
$ws.Range("F44").Value = 258
$ws.Range("F46").Value = 598
$ws.Range("F48").Value = 21

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 22
$ws.Range("F15").Value = 548

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 159

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 159
$ws.Range("F4").Value = 50
$ws.Range("F5").Value = 7121
$ws.Range("F6").Value = 3478
$ws.Range("F8").Value = 3804
$ws.Range("F9").Value = 60
$ws.Range("F11").Value = 66
$ws.Range("F12").Value = 124
$ws.Range("F15").Value = 111
$ws.Range("F16").Value = 350
$ws.Range("F17").Value = 11
$ws.Range("F22").Value = 4055
$ws.Range("F25").Value = 400
$ws.Range("F26").Value = 526
$ws.Range("F27").Value = 1613
$ws.Range("F28").Value = 105
$ws.Range("F29").Value = 89
$ws.Range("F30").Value = 2953
$ws.Range("F31").Value = 2132
$ws.Range("F32").Value = 56
$ws.Range("F33").Value = 75
$ws.Range("F34").Value = 36
$ws.Range("F37").Value = 4153
$ws.Range("F39").Value = 428
$ws.Range("F40").Value = 314
$ws.Range("F41").Value = 22
$ws.Range("F42").Value = 54
$ws.Range("F43").Value = 938
$ws.Range("F44").Value = 739
$ws.Range("F45").Value = 1584
$ws.Range("F46").Value = 258
$ws.Range("F49").Value = 21
